$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates -------------------------------------------------

# Header: "Params" -> "Parameters"
$ws.Range("C1").Value = "Parameters"

# Row 5 (SHOW_ACTIVE_ADMIN_SCREEN): add lobbyActive parameter / comment
$ws.Range("C5").Value = "queue, enqueueKey, lobbyActive "
$ws.Range("D5").Value = "enqueueKey in QR-Code einbinden; lobbyActive for switch-button state"

# Row 7 (SHOW_INVALID_INPUT_ERROR_MESSAGE re-used): now has a Params value
$ws.Range("C7").Value = "msg"

# New "User" group starting at row 9
$ws.Range("A9").Value = "User"
$ws.Range("A9").Font.Bold = $true
$ws.Range("B9").Value = "SHOW_INVALID_INPUT_ERROR_MESSAGE"
$ws.Range("C9").Value = "msg"

$ws.Range("B10").Value = "SHOW_USER_SUCCESSFULLY_ENQUEUED"

# New "Owner" group starting at row 11
$ws.Range("A11").Value = "Owner"
$ws.Range("B11").Value = "UPDATE_OWNER_SCREEEN"
$ws.Range("C11").Value = "queue, enqueueKey, lobbyActive "

# --- Column widths (auto-fit style side effect of the content changes) ---
$ws.Columns.Item(2).ColumnWidth = 37.166666666666664
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668

# --- Selection / view state -----------------------------------------------
$ws.Range("B10").Select() | Out-Null
